$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date and volume updated
$ws.Range("D2").Value = 44291
$ws.Range("J2").Value = 30

# Row 3: date, volume and price fields updated
$ws.Range("D3").Value = 44277
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 11000
$ws.Range("P3").Value = 550

# Row 5: date, volume and price fields updated
$ws.Range("D5").Value = 44284
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 500
